$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings stay as text (preserve exact formatting)
$textCells = @("D2","D3","D4","D5","D6","D7","D9","D11","D12","D13","D14","D15","D16","D17","D18","D19","D21","D22","D23","D25","D26","D28","D29","D30","D31","D32","D34","D35","D36","D37","D38","D40","D43","D44","D45","D46","D48","D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply cell value updates from the diff
$ws.Range('D2').Value = '73.637.29'
$ws.Range('E2').Value = '  +1.13%  '
$ws.Range('D3').Value = '3.985.36'
$ws.Range('E3').Value = '  -1.71%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '615.32'
$ws.Range('E5').Value = '  +10.29%  '
$ws.Range('D6').Value = '169.27'
$ws.Range('E6').Value = '  +11.04%  '
$ws.Range('D7').Value = '0.683'
$ws.Range('E7').Value = '  -1.98%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = '0.768'
$ws.Range('E9').Value = '  +0.95%  '
$ws.Range('E10').Value = '  +8.17%  '
$ws.Range('D11').Value = '55.77'
$ws.Range('E11').Value = '  +3.66%  '
$ws.Range('D12').Value = '0.0000337'
$ws.Range('E12').Value = '  +2.09%  '
$ws.Range('D13').Value = '11.18'
$ws.Range('E13').Value = '  +1.46%  '
$ws.Range('D14').Value = '4.618.17'
$ws.Range('E14').Value = '  -1.59%  '
$ws.Range('D15').Value = '3.974.70'
$ws.Range('E15').Value = '  -1.97%  '
$ws.Range('B16').Value = 'Uniswap'
$ws.Range('C16').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D16').Value = '14.17'
$ws.Range('E16').Value = '  -2.43%  '
$ws.Range('B17').Value = 'Polygon'
$ws.Range('C17').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D17').Value = '1.24'
$ws.Range('E17').Value = '  +2.19%  '
$ws.Range('D18').Value = '20.49'
$ws.Range('E18').Value = '  -1.57%  '
$ws.Range('D19').Value = '73.419.18'
$ws.Range('E19').Value = '  +0.97%  '
$ws.Range('E20').Value = '  -1.05%  '
$ws.Range('D21').Value = '439.96'
$ws.Range('E21').Value = '  -2.07%  '
$ws.Range('D22').Value = '4.87'
$ws.Range('E22').Value = '  +11.18%  '
$ws.Range('D23').Value = '95.98'
$ws.Range('E23').Value = '  -2.19%  '
$ws.Range('E24').Value = '  -5.23%  '
$ws.Range('D25').Value = '14.21'
$ws.Range('E25').Value = '  -3.88%  '
$ws.Range('D26').Value = '4.09'
$ws.Range('E26').Value = '  -3.59%  '
$ws.Range('E27').Value = '  -2.81%  '
$ws.Range('B28').Value = 'Filecoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D28').Value = '10.55'
$ws.Range('E28').Value = '  -3.36%  '
$ws.Range('B29').Value = 'LEO'
$ws.Range('C29').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D29').Value = '5.96'
$ws.Range('E29').Value = '  +0.33%  '
$ws.Range('D30').Value = '36.18'
$ws.Range('E30').Value = '  -3.43%  '
$ws.Range('D31').Value = '7.73'
$ws.Range('E31').Value = '  -1.51%  '
$ws.Range('D32').Value = '13.85'
$ws.Range('E32').Value = '  +1.08%  '
$ws.Range('E33').Value = '  +20.33%  '
$ws.Range('D34').Value = '0.130'
$ws.Range('E34').Value = '  -3.87%  '
$ws.Range('D35').Value = '47.73'
$ws.Range('E35').Value = '  -2.20%  '
$ws.Range('D36').Value = '70.99'
$ws.Range('E36').Value = '  +5.44%  '
$ws.Range('D37').Value = '648.23'
$ws.Range('E37').Value = '  -7.11%  '
$ws.Range('D38').Value = '0.430'
$ws.Range('E38').Value = '  -4.87%  '
$ws.Range('E39').Value = '  +0.19%  '
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  -0.12%  '
$ws.Range('E41').Value = '  -2.62%  '
$ws.Range('E42').Value = '  +0.11%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = '0.0486'
$ws.Range('E43').Value = '  -2.40%  '
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').Value = '3.24'
$ws.Range('E44').Value = '  +43.71%  '
$ws.Range('B45').Value = 'WEMIXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D45').Value = '3.22'
$ws.Range('E45').Value = '  -6.70%  '
$ws.Range('B46').Value = 'THORChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D46').Value = '10.62'
$ws.Range('E46').Value = '  -5.80%  '
$ws.Range('E47').Value = '  -2.07%  '
$ws.Range('D48').Value = '0.000301'
$ws.Range('E48').Value = '  +8.08%  '
$ws.Range('D49').Value = '3.42'
$ws.Range('E49').Value = '  +2.20%  '
$ws.Range('E50').Value = '  -5.20%  '
$ws.Range('E51').Value = '  -4.31%  '
